$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextCell 2 1 "2026-02-01"
$ws.Cells.Item(2,2).Value = "기술"
$ws.Cells.Item(2,3).Value = "AI 훈풍 탄 반도체…1월 수출 역대 최고, 흑자 행진"
Set-TextCell 2 4 "2026-02-01"
$ws.Cells.Item(2,5).Value = "https://n.news.naver.com/mnews/article/005/0001829492?sid=101"

# Row 3
Set-TextCell 3 1 "2026-02-01"
$ws.Cells.Item(3,2).Value = "기업"
$ws.Cells.Item(3,3).Value = "젠슨 황, '오픈AI에 대한 불만설' 부인…""최대 규모 투자할 것"""
Set-TextCell 3 4 "2026-02-01"
$ws.Cells.Item(3,5).Value = "https://n.news.naver.com/mnews/article/001/0015877800?sid=104"

# Row 4
Set-TextCell 4 1 "2026-02-01"
$ws.Cells.Item(4,2).Value = "기업"
$ws.Cells.Item(4,3).Value = "“오픈AI·엔비디아 145조원 초대형 계약 ‘제동’”…AI동맹 흔들리나"
Set-TextCell 4 4 "2026-02-01"
$ws.Cells.Item(4,5).Value = "https://n.news.naver.com/mnews/article/018/0006210754?sid=101"

# Row 5
Set-TextCell 5 1 "2026-02-01"
$ws.Cells.Item(5,2).Value = "기술"
$ws.Cells.Item(5,3).Value = "정부, 독자AI 모델 기반 '국방 AI' 개발한다"
Set-TextCell 5 4 "2026-01-30"
$ws.Cells.Item(5,5).Value = "https://n.news.naver.com/mnews/article/003/0013740053?sid=105"

# Row 6
Set-TextCell 6 1 "2026-02-01"
$ws.Cells.Item(6,2).Value = "산업"
$ws.Cells.Item(6,3).Value = "인천항만공사, AI 기반 혁신 우수사례 발굴 주력"
Set-TextCell 6 4 "2026-02-01"
$ws.Cells.Item(6,5).Value = "https://n.news.naver.com/mnews/article/005/0001829508?sid=102"

# Row 7
Set-TextCell 7 1 "2026-02-01"
$ws.Cells.Item(7,2).Value = "정책"
$ws.Cells.Item(7,3).Value = "美 'AI규제완화' 슈퍼팩, 중간선거 앞두고 1천800억원 모금"
Set-TextCell 7 4 "2026-02-01"
$ws.Cells.Item(7,5).Value = "https://n.news.naver.com/mnews/article/001/0015878085?sid=104"

# Row 8
Set-TextCell 8 1 "2026-02-01"
$ws.Cells.Item(8,2).Value = "정책"
$ws.Cells.Item(8,3).Value = "미 'AI규제완화' 슈퍼팩, 중간선거 앞두고 1천800억 원 모금"
Set-TextCell 8 4 "2026-02-01"
$ws.Cells.Item(8,5).Value = "https://n.news.naver.com/mnews/article/055/0001329287?sid=104"

# Row 9
Set-TextCell 9 1 "2026-02-01"
$ws.Cells.Item(9,2).Value = "산업"
$ws.Cells.Item(9,3).Value = "[인간과 AI] 의료 현장까지 파고든 AI‥사고시 책임은 누구에게?"
Set-TextCell 9 4 "2026-01-31"
$ws.Cells.Item(9,5).Value = "https://n.news.naver.com/mnews/article/214/0001477886?sid=102"

# Row 10
Set-TextCell 10 1 "2026-02-01"
$ws.Cells.Item(10,2).Value = "정부(과기부)"
$ws.Cells.Item(10,3).Value = "오늘의인사-헌법재판소, 대법원, 국회, 감사원 외"
Set-TextCell 10 4 "2026-02-01"
$ws.Cells.Item(10,5).Value = "https://n.news.naver.com/mnews/article/032/0003425035?sid=102"

# Row 11
Set-TextCell 11 1 "2026-02-01"
$ws.Cells.Item(11,2).Value = "정부(과기부)"
$ws.Cells.Item(11,3).Value = "국방 안보, 'K-AI'로 무장…민·관·군 AI 대전환 본격화"
Set-TextCell 11 4 "2026-02-01"
$ws.Cells.Item(11,5).Value = "https://n.news.naver.com/mnews/article/092/0002408658?sid=105"

# Row 12
Set-TextCell 12 1 "2026-02-01"
$ws.Cells.Item(12,2).Value = "정부(과기부)"
$ws.Cells.Item(12,3).Value = "KAIST, 모자처럼 쓰는 탈모 치료기 개발… 光치료 기술, 탈모 치료의 패..."
Set-TextCell 12 4 "2026-02-01"
$ws.Cells.Item(12,5).Value = "http://www.biotimes.co.kr/news/articleView.html?idxno=26823"

# Row 13
Set-TextCell 13 1 "2026-02-01"
$ws.Cells.Item(13,2).Value = "정부(과기부)"
$ws.Cells.Item(13,3).Value = "건설연, 국가연구개발 우수성과 100선 선정…“스마트건설기술 상용화”"
Set-TextCell 13 4 "2026-02-01"
$ws.Cells.Item(13,5).Value = "https://n.news.naver.com/mnews/article/666/0000094935?sid=101"

# Row 14 (new)
Set-TextCell 14 1 "2026-02-01"
$ws.Cells.Item(14,2).Value = "정부(과기부)"
$ws.Cells.Item(14,3).Value = "KAIST, 헬맷은 가라...모자처럼 쓰는 탈모 예방 OLED 개발"
Set-TextCell 14 4 "2026-02-01"
$ws.Cells.Item(14,5).Value = "http://www.veritas-a.com/news/articleView.html?idxno=596264"
